$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

# Column A holds a date-like text value ("2025-10-05"). A plain assignment
# would be auto-parsed into a date serial number by Excel, which does not
# match the source data (the workbook stores these as literal text). Force
# text entry with a leading apostrophe, then clear the resulting cell style
# back to Normal so no stray formatting is left behind.
$cellA = $ws.Cells.Item($row, 1)
$cellA.Value = "'2025-10-05"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "21:19:13"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,794.1737"
